$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update existing cell values ---
$ws.Range("A120").Value = 6110051104

$ws.Range("B227").Value = 6110051103
$ws.Range("A228").Value = 6110051103
$ws.Range("A229").Value = 6110051103
$ws.Range("B230").Value = 6110051103

# --- Append new row 231 ---
$ws.Range("A231").Value = 6110051001
$ws.Range("B231").Value = 6110051104
$ws.Range("C231").Value = 2.1

# --- Update selection / active cell to match the new data extent ---
$null = $ws.Range("C2:C231").Select()

Write-Host "done"
